$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 646 (the "Super star" song entry) — this shifts all following
# rows up by one and removes the now-unused shared string automatically.
$ws.Rows.Item(646).Delete()

# Re-create the checkbox data validation on column E so it reflects the new
# row count / boundary around the deleted row.
$ws.Cells.Validation.Delete()
$r1 = $ws.Range("E2:E645")
$r2 = $ws.Range("E646:E1048576")
$u = $excel.Union($r1, $r2)
$u.Validation.Add(7, 1, 1, "IF(TRUE,OR(E2=0,E2=1),""Checkbox"")")
$u.Validation.ErrorTitle = "输入内容有误"
$u.Validation.ErrorMessage = "请选择勾选或取消勾选"
$u.Validation.IgnoreBlank = $true
$u.Validation.ShowError = $true
$u.Validation.ShowInput = $false

# Restore the frozen-pane view and put the active cell where the user left
# off (on the row that used to be 647, now 646) after deleting the row.
$ws.Range("A646").Select()
